$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting existing data rows down
$ws.Rows.Item(2).Insert()

# The inserted row picks up header formatting by default; clear it so it
# matches the plain formatting used by the other data rows
$ws.Range("A2:R2").ClearFormats()

# Re-apply the date number format used by column D in the other data rows
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new first data row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44638
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112030
$ws.Range("G2").Value = "Poroto granado"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2650
$ws.Range("N2").Value = "`$/kilo"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 2650
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
